# Apply the update described by the commit diff:
#  1. Update the "Förändrad" (Changed) date in column C for all existing
#     data rows (rows 2-406) from 45175 (2023-09-06) to 45177 (2023-09-08).
#  2. Give row 406 an explicit custom row height (15), matching the diff.
#  3. Append a brand-new data row (row 407) for case "A 41441-2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bulk-update column C ("Förändrad") for rows 2 through 406.
$ws.Range("C2:C406").Value = 45177

# 2. Row 406 becomes an explicit custom-height row (15pt).
$ws.Rows.Item(406).RowHeight = 15

# 3. Append the new row (407) with the new cleaving/felling notice.
$ws.Range("A407").Value = "A 41441-2023"

$ws.Range("B407").Value = 45175
$ws.Range("C407").Value = 45177
$ws.Range("B407:C407").NumberFormat = "YYYY-MM-DD"

$ws.Range("D407").Value = "SÖDERMANLANDS LÄN"
$ws.Range("E407").Value = "VINGÅKER"

$ws.Range("G407").Value = 10.9
$ws.Range("H407").Value = 0
$ws.Range("I407").Value = 0
$ws.Range("J407").Value = 0
$ws.Range("K407").Value = 0
$ws.Range("L407").Value = 0
$ws.Range("M407").Value = 0
$ws.Range("N407").Value = 0
$ws.Range("O407").Value = 0
$ws.Range("P407").Value = 0
$ws.Range("Q407").Value = 0

# R407 keeps the same wrap-text style used by the rest of the "Artnamn"
# column, but with no content (same as R406/R405/...).
$ws.Range("R407").WrapText = $true
